# Auto-generated script to apply two-digit multiplication worksheet updates
$d = $word.ActiveDocument

$pairs = @(
    @("2023-05-25 Thursday", "2023-05-26 Friday"),
    @("59×43=2537", "61×81=4941"),
    @("45×10=450", "77×67=5159"),
    @("74×39=2886", "29×87=2523"),
    @("79×63=4977", "13×20=260"),
    @("36×27=972", "43×77=3311"),
    @("53×96=5088", "100×22=2200"),
    @("53×81=4293", "87×10=870"),
    @("45×59=2655", "43×71=3053"),
    @("68×52=3536", "77×32=2464"),
    @("70×82=5740", "65×92=5980"),
    @("22×24=528", "85×94=7990"),
    @("99×36=3564", "63×75=4725"),
    @("15×15=225", "76×61=4636"),
    @("37×93=3441", "98×96=9408"),
    @("30×92=2760", "81×28=2268"),
    @("20×38=760", "20×46=920"),
    @("20×56=1120", "57×77=4389"),
    @("88×77=6776", "94×11=1034"),
    @("78×65=5070", "84×97=8148"),
    @("17×58=986", "83×95=7885"),
    @("19×24=456", "26×29=754"),
    @("70×30=2100", "74×45=3330"),
    @("18×70=1260", "94×19=1786"),
    @("73×78=5694", "35×75=2625"),
    @("33×49=1617", "44×95=4180"),
    @("40×26=1040", "49×78=3822"),
    @("62×83=5146", "68×66=4488"),
    @("26×72=1872", "46×55=2530"),
    @("15×13=195", "99×27=2673"),
    @("90×17=1530", "25×32=800"),
    @("50×94=4700", "95×87=8265"),
    @("61×33=2013", "33×11=363"),
    @("22×54=1188", "62×68=4216"),
    @("93×54=5022", "91×16=1456"),
    @("76×33=2508", "46×10=460"),
    @("31×43=1333", "52×61=3172"),
    @("29×15=435", "98×80=7840"),
    @("42×90=3780", "12×11=132"),
    @("69×94=6486", "90×34=3060"),
    @("59×60=3540", "12×84=1008"),
    @("34×30=1020", "30×39=1170"),
    @("90×28=2520", "71×86=6106"),
    @("13×57=741", "36×96=3456"),
    @("14×34=476", "39×76=2964"),
    @("66×30=1980", "97×60=5820"),
    @("65×37=2405", "19×43=817"),
    @("95×74=7030", "100×42=4200"),
    @("62×22=1364", "43×67=2881"),
    @("30×65=1950", "90×27=2430"),
    @("88×44=3872", "90×95=8550"),
    @("56×69=3864", "14×63=882"),
    @("77×71=5467", "20×89=1780"),
    @("90×86=7740", "13×78=1014"),
    @("47×86=4042", "97×40=3880"),
    @("39×78=3042", "34×74=2516"),
    @("79×83=6557", "26×20=520"),
    @("68×78=5304", "30×88=2640"),
    @("40×17=680", "98×98=9604"),
    @("40×30=1200", "27×86=2322"),
    @("53×53=2809", "73×33=2409"),
    @("66×20=1320", "54×10=540"),
    @("70×46=3220", "68×17=1156"),
    @("27×17=459", "28×27=756"),
    @("53×14=742", "70×53=3710"),
    @("96×45=4320", "11×53=583"),
    @("81×85=6885", "35×15=525"),
    @("60×60=3600", "82×95=7790"),
    @("55×20=1100", "76×72=5472"),
    @("82×69=5658", "29×24=696"),
    @("30×16=480", "98×40=3920"),
    @("33×96=3168", "11×72=792"),
    @("33×82=2706", "11×49=539"),
    @("66×58=3828", "68×67=4556"),
    @("63×63=3969", "38×51=1938"),
    @("14×23=322", "84×85=7140"),
    @("87×27=2349", "81×38=3078"),
    @("66×23=1518", "85×95=8075"),
    @("22×99=2178", "30×97=2910"),
    @("35×70=2450", "33×99=3267"),
    @("76×83=6308", "29×28=812"),
    @("92×80=7360", "65×63=4095"),
    @("39×77=3003", "11×53=583"),
    @("40×45=1800", "33×75=2475"),
    @("39×59=2301", "72×13=936"),
    @("89×31=2759", "86×37=3182"),
    @("91×34=3094", "42×27=1134"),
    @("75×62=4650", "60×67=4020"),
    @("24×50=1200", "99×54=5346"),
    @("18×97=1746", "82×74=6068"),
    @("46×54=2484", "87×29=2523"),
    @("51×53=2703", "45×97=4365"),
    @("42×28=1176", "87×51=4437"),
    @("34×92=3128", "52×92=4784"),
    @("37×100=3700", "23×41=943"),
    @("29×68=1972", "14×77=1078"),
    @("12×15=180", "94×58=5452"),
    @("45×44=1980", "80×79=6320"),
    @("51×18=918", "30×42=1260"),
    @("60×32=1920", "15×93=1395"),
    @("12×34=408", "84×44=3696"),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "Applied $($pairs.Count) replacements"
